$wb = $excel.ActiveWorkbook

$headers = @('Turn','Ticker','PriceVol_Signal','PriceVol_Reason','MACD_Signal','MACD_Reason','Bollinger_Signal','Bollinger_Reason','Volume_Signal','Volume_Reason','Overall_Signal','Signal_Reasoning')

# ---- Sheet: NFLX ----
$ws = $wb.Worksheets.Item('NFLX')

# Insert new columns to grow from 6 to 12 columns, preserving existing data/style
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(8).Insert()
$ws.Columns.Item(10).Insert()
$ws.Columns.Item(11).Insert()
$ws.Columns.Item(12).Insert()

# Write header row
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Row 2
$ws.Cells.Item(2, 3).Value = 'buy'
$ws.Cells.Item(2, 4).Value = 'Trend: uptrend, Volatility: 15.02854220982136'
$ws.Cells.Item(2, 5).Value = 'sell'
$ws.Cells.Item(2, 6).Value = 'Trend: Neutral, Histogram Strength: -1.95'
$ws.Cells.Item(2, 7).Value = 'Sell'
$ws.Cells.Item(2, 8).Value = 'Band: neutral, Volatility: 15.02854220982136'
$ws.Cells.Item(2, 9).Value = 'Sell'
$ws.Cells.Item(2, 10).Value = 'Trend: downtrend'
$ws.Cells.Item(2, 11).Value = 'Sell'
$ws.Cells.Item(2, 12).Value = 'PriceVol: buy (+0.40) | MACD: sell (-0.25) | Bollinger: sell (-0.20) | Volume: sell (-0.15)'
# Row 3
$ws.Cells.Item(3, 3).Value = 'sell'
$ws.Cells.Item(3, 4).Value = 'Trend: downtrend, Volatility: 15.47683199402998'
$ws.Cells.Item(3, 5).Value = 'sell'
$ws.Cells.Item(3, 6).Value = 'Trend: Neutral, Histogram Strength: -3.97'
$ws.Cells.Item(3, 7).Value = 'Sell'
$ws.Cells.Item(3, 8).Value = 'Band: neutral, Volatility: 15.47683199402998'
$ws.Cells.Item(3, 9).Value = 'Sell'
$ws.Cells.Item(3, 10).Value = 'Trend: downtrend'
$ws.Cells.Item(3, 11).Value = 'Sell'
$ws.Cells.Item(3, 12).Value = 'PriceVol: sell (-0.40) | MACD: sell (-0.25) | Bollinger: sell (-0.20) | Volume: sell (-0.15)'
# Row 4
$ws.Cells.Item(4, 3).Value = 'sell'
$ws.Cells.Item(4, 4).Value = 'Trend: downtrend, Volatility: 25.65931112908225'
$ws.Cells.Item(4, 5).Value = 'sell'
$ws.Cells.Item(4, 6).Value = 'Trend: Neutral, Histogram Strength: -8.8'
$ws.Cells.Item(4, 7).Value = 'Buy'
$ws.Cells.Item(4, 8).Value = 'Band: oversold, Volatility: 25.65931112908225'
$ws.Cells.Item(4, 9).Value = 'Sell'
$ws.Cells.Item(4, 10).Value = 'Trend: downtrend'
$ws.Cells.Item(4, 11).Value = 'Sell'
$ws.Cells.Item(4, 12).Value = 'PriceVol: sell (-0.40) | MACD: sell (-0.25) | Bollinger: buy (+0.20) | Volume: sell (-0.15)'
# Row 5
$ws.Cells.Item(5, 3).Value = 'buy'
$ws.Cells.Item(5, 4).Value = 'Trend: uptrend, Volatility: 11.13681988463755'
$ws.Cells.Item(5, 5).Value = 'sell'
$ws.Cells.Item(5, 6).Value = 'Trend: Sell, Histogram Strength: -4.05'
$ws.Cells.Item(5, 7).Value = 'Sell'
$ws.Cells.Item(5, 8).Value = 'Band: neutral, Volatility: 11.13681988463755'
$ws.Cells.Item(5, 9).Value = 'Buy'
$ws.Cells.Item(5, 10).Value = 'Trend: uptrend'
$ws.Cells.Item(5, 11).Value = 'Buy'
$ws.Cells.Item(5, 12).Value = 'PriceVol: buy (+0.40) | MACD: sell (-0.25) | Bollinger: sell (-0.20) | Volume: buy (+0.15)'
# Row 6
$ws.Cells.Item(6, 3).Value = 'sell'
$ws.Cells.Item(6, 4).Value = 'Trend: downtrend, Volatility: 11.88190320653737'
$ws.Cells.Item(6, 5).Value = 'sell'
$ws.Cells.Item(6, 6).Value = 'Trend: Sell, Histogram Strength: 1.83'
$ws.Cells.Item(6, 7).Value = 'Sell'
$ws.Cells.Item(6, 8).Value = 'Band: neutral, Volatility: 11.88190320653737'
$ws.Cells.Item(6, 9).Value = 'Buy'
$ws.Cells.Item(6, 10).Value = 'Trend: uptrend'
$ws.Cells.Item(6, 11).Value = 'Sell'
$ws.Cells.Item(6, 12).Value = 'PriceVol: sell (-0.40) | MACD: sell (-0.25) | Bollinger: sell (-0.20) | Volume: buy (+0.15)'

# ---- Sheet: PG ----
$ws = $wb.Worksheets.Item('PG')

# Insert new columns to grow from 6 to 12 columns, preserving existing data/style
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(8).Insert()
$ws.Columns.Item(10).Insert()
$ws.Columns.Item(11).Insert()
$ws.Columns.Item(12).Insert()

# Write header row
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Row 2
$ws.Cells.Item(2, 3).Value = 'buy'
$ws.Cells.Item(2, 4).Value = 'Trend: uptrend, Volatility: 0.840093532207674'
$ws.Cells.Item(2, 5).Value = 'buy'
$ws.Cells.Item(2, 6).Value = 'Trend: Buy, Histogram Strength: 0.39'
$ws.Cells.Item(2, 7).Value = 'Buy'
$ws.Cells.Item(2, 8).Value = 'Band: neutral, Volatility: 0.840093532207674'
$ws.Cells.Item(2, 9).Value = 'Sell'
$ws.Cells.Item(2, 10).Value = 'Trend: downtrend'
$ws.Cells.Item(2, 11).Value = 'Buy'
$ws.Cells.Item(2, 12).Value = 'PriceVol: buy (+0.40) | MACD: buy (+0.25) | Bollinger: buy (+0.20) | Volume: sell (-0.15)'
# Row 3
$ws.Cells.Item(3, 3).Value = 'sell'
$ws.Cells.Item(3, 4).Value = 'Trend: downtrend, Volatility: 0.9731857357105586'
$ws.Cells.Item(3, 5).Value = 'buy'
$ws.Cells.Item(3, 6).Value = 'Trend: Buy, Histogram Strength: 0.21'
$ws.Cells.Item(3, 7).Value = 'Buy'
$ws.Cells.Item(3, 8).Value = 'Band: neutral, Volatility: 0.9731857357105586'
$ws.Cells.Item(3, 9).Value = 'Buy'
$ws.Cells.Item(3, 10).Value = 'Trend: uptrend'
$ws.Cells.Item(3, 11).Value = 'Buy'
$ws.Cells.Item(3, 12).Value = 'PriceVol: sell (-0.40) | MACD: buy (+0.25) | Bollinger: buy (+0.20) | Volume: buy (+0.15)'
# Row 4
$ws.Cells.Item(4, 3).Value = 'buy'
$ws.Cells.Item(4, 4).Value = 'Trend: uptrend, Volatility: 2.476546174086182'
$ws.Cells.Item(4, 5).Value = 'sell'
$ws.Cells.Item(4, 6).Value = 'Trend: Neutral, Histogram Strength: 0.19'
$ws.Cells.Item(4, 7).Value = 'Sell'
$ws.Cells.Item(4, 8).Value = 'Band: overbought, Volatility: 2.476546174086182'
$ws.Cells.Item(4, 9).Value = 'Buy'
$ws.Cells.Item(4, 10).Value = 'Trend: uptrend'
$ws.Cells.Item(4, 11).Value = 'Buy'
$ws.Cells.Item(4, 12).Value = 'PriceVol: buy (+0.40) | MACD: sell (-0.25) | Bollinger: sell (-0.20) | Volume: buy (+0.15)'
# Row 5
$ws.Cells.Item(5, 3).Value = 'buy'
$ws.Cells.Item(5, 4).Value = 'Trend: uptrend, Volatility: 2.935980407356449'
$ws.Cells.Item(5, 5).Value = 'buy'
$ws.Cells.Item(5, 6).Value = 'Trend: Buy, Histogram Strength: 0.93'
$ws.Cells.Item(5, 7).Value = 'Sell'
$ws.Cells.Item(5, 8).Value = 'Band: overbought, Volatility: 2.935980407356449'
$ws.Cells.Item(5, 9).Value = 'Buy'
$ws.Cells.Item(5, 10).Value = 'Trend: uptrend'
$ws.Cells.Item(5, 11).Value = 'Buy'
$ws.Cells.Item(5, 12).Value = 'PriceVol: buy (+0.40) | MACD: buy (+0.25) | Bollinger: sell (-0.20) | Volume: buy (+0.15)'
# Row 6
$ws.Cells.Item(6, 3).Value = 'buy'
$ws.Cells.Item(6, 4).Value = 'Trend: uptrend, Volatility: 1.772159027584043'
$ws.Cells.Item(6, 5).Value = 'buy'
$ws.Cells.Item(6, 6).Value = 'Trend: Buy, Histogram Strength: 0.34'
$ws.Cells.Item(6, 7).Value = 'Buy'
$ws.Cells.Item(6, 8).Value = 'Band: neutral, Volatility: 1.772159027584043'
$ws.Cells.Item(6, 9).Value = 'Sell'
$ws.Cells.Item(6, 10).Value = 'Trend: downtrend'
$ws.Cells.Item(6, 11).Value = 'Buy'
$ws.Cells.Item(6, 12).Value = 'PriceVol: buy (+0.40) | MACD: buy (+0.25) | Bollinger: buy (+0.20) | Volume: sell (-0.15)'

# ---- Sheet: TSLA ----
$ws = $wb.Worksheets.Item('TSLA')

# Insert new columns to grow from 6 to 12 columns, preserving existing data/style
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(8).Insert()
$ws.Columns.Item(10).Insert()
$ws.Columns.Item(11).Insert()
$ws.Columns.Item(12).Insert()

# Write header row
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Row 2
$ws.Cells.Item(2, 3).Value = 'sell'
$ws.Cells.Item(2, 4).Value = 'Trend: downtrend, Volatility: 26.23960165928774'
$ws.Cells.Item(2, 5).Value = 'sell'
$ws.Cells.Item(2, 6).Value = 'Trend: Neutral, Histogram Strength: -7.68'
$ws.Cells.Item(2, 7).Value = 'Sell'
$ws.Cells.Item(2, 8).Value = 'Band: neutral, Volatility: 26.23960165928774'
$ws.Cells.Item(2, 9).Value = 'Buy'
$ws.Cells.Item(2, 10).Value = 'Trend: uptrend'
$ws.Cells.Item(2, 11).Value = 'Sell'
$ws.Cells.Item(2, 12).Value = 'PriceVol: sell (-0.40) | MACD: sell (-0.25) | Bollinger: sell (-0.20) | Volume: buy (+0.15)'
# Row 3
$ws.Cells.Item(3, 3).Value = 'buy'
$ws.Cells.Item(3, 4).Value = 'Trend: uptrend, Volatility: 11.07897708876292'
$ws.Cells.Item(3, 5).Value = 'sell'
$ws.Cells.Item(3, 6).Value = 'Trend: Neutral, Histogram Strength: -2.91'
$ws.Cells.Item(3, 7).Value = 'Sell'
$ws.Cells.Item(3, 8).Value = 'Band: neutral, Volatility: 11.07897708876292'
$ws.Cells.Item(3, 9).Value = 'Buy'
$ws.Cells.Item(3, 10).Value = 'Trend: uptrend'
$ws.Cells.Item(3, 11).Value = 'Buy'
$ws.Cells.Item(3, 12).Value = 'PriceVol: buy (+0.40) | MACD: sell (-0.25) | Bollinger: sell (-0.20) | Volume: buy (+0.15)'
# Row 4
$ws.Cells.Item(4, 3).Value = 'sell'
$ws.Cells.Item(4, 4).Value = 'Trend: downtrend, Volatility: 17.71938661197847'
$ws.Cells.Item(4, 5).Value = 'sell'
$ws.Cells.Item(4, 6).Value = 'Trend: Neutral, Histogram Strength: -6.47'
$ws.Cells.Item(4, 7).Value = 'Sell'
$ws.Cells.Item(4, 8).Value = 'Band: neutral, Volatility: 17.71938661197847'
$ws.Cells.Item(4, 9).Value = 'Buy'
$ws.Cells.Item(4, 10).Value = 'Trend: uptrend'
$ws.Cells.Item(4, 11).Value = 'Sell'
$ws.Cells.Item(4, 12).Value = 'PriceVol: sell (-0.40) | MACD: sell (-0.25) | Bollinger: sell (-0.20) | Volume: buy (+0.15)'
# Row 5
$ws.Cells.Item(5, 3).Value = 'buy'
$ws.Cells.Item(5, 4).Value = 'Trend: uptrend, Volatility: 14.21137321202383'
$ws.Cells.Item(5, 5).Value = 'sell'
$ws.Cells.Item(5, 6).Value = 'Trend: Neutral, Histogram Strength: -6.27'
$ws.Cells.Item(5, 7).Value = 'Sell'
$ws.Cells.Item(5, 8).Value = 'Band: neutral, Volatility: 14.21137321202383'
$ws.Cells.Item(5, 9).Value = 'Buy'
$ws.Cells.Item(5, 10).Value = 'Trend: uptrend'
$ws.Cells.Item(5, 11).Value = 'Buy'
$ws.Cells.Item(5, 12).Value = 'PriceVol: buy (+0.40) | MACD: sell (-0.25) | Bollinger: sell (-0.20) | Volume: buy (+0.15)'
# Row 6
$ws.Cells.Item(6, 3).Value = 'buy'
$ws.Cells.Item(6, 4).Value = 'Trend: uptrend, Volatility: 25.30465155808687'
$ws.Cells.Item(6, 5).Value = 'sell'
$ws.Cells.Item(6, 6).Value = 'Trend: Sell, Histogram Strength: 2.52'
$ws.Cells.Item(6, 7).Value = 'Sell'
$ws.Cells.Item(6, 8).Value = 'Band: neutral, Volatility: 25.30465155808687'
$ws.Cells.Item(6, 9).Value = 'Buy'
$ws.Cells.Item(6, 10).Value = 'Trend: uptrend'
$ws.Cells.Item(6, 11).Value = 'Buy'
$ws.Cells.Item(6, 12).Value = 'PriceVol: buy (+0.40) | MACD: sell (-0.25) | Bollinger: sell (-0.20) | Volume: buy (+0.15)'

# ---- Sheet: XOM ----
$ws = $wb.Worksheets.Item('XOM')

# Insert new columns to grow from 6 to 12 columns, preserving existing data/style
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(8).Insert()
$ws.Columns.Item(10).Insert()
$ws.Columns.Item(11).Insert()
$ws.Columns.Item(12).Insert()

# Write header row
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Row 2
$ws.Cells.Item(2, 3).Value = 'buy'
$ws.Cells.Item(2, 4).Value = 'Trend: uptrend, Volatility: 0.920212190850825'
$ws.Cells.Item(2, 5).Value = 'sell'
$ws.Cells.Item(2, 6).Value = 'Trend: Neutral, Histogram Strength: -0.31'
$ws.Cells.Item(2, 7).Value = 'Buy'
$ws.Cells.Item(2, 8).Value = 'Band: neutral, Volatility: 0.920212190850825'
$ws.Cells.Item(2, 9).Value = 'Buy'
$ws.Cells.Item(2, 10).Value = 'Trend: uptrend'
$ws.Cells.Item(2, 11).Value = 'Buy'
$ws.Cells.Item(2, 12).Value = 'PriceVol: buy (+0.40) | MACD: sell (-0.25) | Bollinger: buy (+0.20) | Volume: buy (+0.15)'
# Row 3
$ws.Cells.Item(3, 3).Value = 'buy'
$ws.Cells.Item(3, 4).Value = 'Trend: uptrend, Volatility: 1.526912539984496'
$ws.Cells.Item(3, 5).Value = 'sell'
$ws.Cells.Item(3, 6).Value = 'Trend: Neutral, Histogram Strength: -0.41'
$ws.Cells.Item(3, 7).Value = 'Buy'
$ws.Cells.Item(3, 8).Value = 'Band: neutral, Volatility: 1.526912539984496'
$ws.Cells.Item(3, 9).Value = 'Sell'
$ws.Cells.Item(3, 10).Value = 'Trend: downtrend'
$ws.Cells.Item(3, 11).Value = 'Buy'
$ws.Cells.Item(3, 12).Value = 'PriceVol: buy (+0.40) | MACD: sell (-0.25) | Bollinger: buy (+0.20) | Volume: sell (-0.15)'
# Row 4
$ws.Cells.Item(4, 3).Value = 'sell'
$ws.Cells.Item(4, 4).Value = 'Trend: downtrend, Volatility: 0.7703122434318123'
$ws.Cells.Item(4, 5).Value = 'sell'
$ws.Cells.Item(4, 6).Value = 'Trend: Neutral, Histogram Strength: -0.3'
$ws.Cells.Item(4, 7).Value = 'Buy'
$ws.Cells.Item(4, 8).Value = 'Band: neutral, Volatility: 0.7703122434318123'
$ws.Cells.Item(4, 9).Value = 'Sell'
$ws.Cells.Item(4, 10).Value = 'Trend: downtrend'
$ws.Cells.Item(4, 11).Value = 'Sell'
$ws.Cells.Item(4, 12).Value = 'PriceVol: sell (-0.40) | MACD: sell (-0.25) | Bollinger: buy (+0.20) | Volume: sell (-0.15)'
# Row 5
$ws.Cells.Item(5, 3).Value = 'buy'
$ws.Cells.Item(5, 4).Value = 'Trend: uptrend, Volatility: 0.6405838111032902'
$ws.Cells.Item(5, 5).Value = 'sell'
$ws.Cells.Item(5, 6).Value = 'Trend: Neutral, Histogram Strength: -0.01'
$ws.Cells.Item(5, 7).Value = 'Buy'
$ws.Cells.Item(5, 8).Value = 'Band: neutral, Volatility: 0.6405838111032902'
$ws.Cells.Item(5, 9).Value = 'Sell'
$ws.Cells.Item(5, 10).Value = 'Trend: downtrend'
$ws.Cells.Item(5, 11).Value = 'Buy'
$ws.Cells.Item(5, 12).Value = 'PriceVol: buy (+0.40) | MACD: sell (-0.25) | Bollinger: buy (+0.20) | Volume: sell (-0.15)'
# Row 6
$ws.Cells.Item(6, 3).Value = 'buy'
$ws.Cells.Item(6, 4).Value = 'Trend: uptrend, Volatility: 0.9047546365622828'
$ws.Cells.Item(6, 5).Value = 'sell'
$ws.Cells.Item(6, 6).Value = 'Trend: Sell, Histogram Strength: 0.04'
$ws.Cells.Item(6, 7).Value = 'Buy'
$ws.Cells.Item(6, 8).Value = 'Band: neutral, Volatility: 0.9047546365622828'
$ws.Cells.Item(6, 9).Value = 'Sell'
$ws.Cells.Item(6, 10).Value = 'Trend: downtrend'
$ws.Cells.Item(6, 11).Value = 'Buy'
$ws.Cells.Item(6, 12).Value = 'PriceVol: buy (+0.40) | MACD: sell (-0.25) | Bollinger: buy (+0.20) | Volume: sell (-0.15)'
